$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.399.31'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.09%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.976.28'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.14%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -11.63%  '
$ws.Range("E6").Value = '  -3.42%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.36'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.13%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '58.61'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.07%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.370'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0747'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.14%  '
$ws.Range("E12").Value = '  -3.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.268.00'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '13.88'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.92'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.748'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -7.23%  '
$ws.Range("E17").Value = '  -4.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.970.61'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.314.91'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '67.54'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0804'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.88%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '221.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.40%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -11.91%  '
$ws.Range("E27").Value = '  -2.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.15%  '
$ws.Range("E30").Value = '  -3.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.32'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("E32").Value = '  -3.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.34'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0604'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.48%  '
$ws.Range("E35").Value = '  -6.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.29'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.34%  '
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("E38").Value = '  -3.14%  '
$ws.Range("E39").Value = '  -3.40%  '
$ws.Range("E40").Value = '  +3.04%  '
$ws.Range("E41").Value = '  -1.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.450.12'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0897'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0200'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.01%  '
$ws.Range("E45").Value = '  -10.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '87.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.37%  '
$ws.Range("E47").Value = '  -3.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '14.77'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.00%  '
$ws.Range("E49").Value = '  -1.87%  '
$ws.Range("E50").Value = '  -4.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.160.78'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.04%  '
